# Updated the supplement PDUFA fees based on year
# Adds a new "SupplementFees" column (C) to Table1, mirroring the existing
# "Fees" column's currency formatting for matching rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the existing table (Table1, currently A1:B33) with a new column.
$lo = $ws.ListObjects.Item(1)
$col = $lo.ListColumns.Add()

# Header cell (shared string "SupplementFees"); reuse the header style from B1.
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("C1").Value = "SupplementFees"

# Data rows: (row, style-source cell to copy number format from, value)
$data = @(
    @{ Row = 2;  StyleFrom = "B7";  Value = 50000 },
    @{ Row = 3;  StyleFrom = "B7";  Value = 75000 },
    @{ Row = 4;  StyleFrom = "B20"; Value = 104000 },
    @{ Row = 5;  StyleFrom = "B20"; Value = 102000 },
    @{ Row = 6;  StyleFrom = "B20"; Value = 102500 },
    @{ Row = 7;  StyleFrom = "B7";  Value = 128423 },
    @{ Row = 8;  StyleFrom = "B7";  Value = 136141 },
    @{ Row = 9;  StyleFrom = "B20"; Value = 142870 },
    @{ Row = 10; StyleFrom = "B7";  Value = 133803 },
    @{ Row = 11; StyleFrom = "B7";  Value = 129226 },
    @{ Row = 12; StyleFrom = "B7";  Value = 266700 },
    @{ Row = 13; StyleFrom = "B7";  Value = 286750 },
    @{ Row = 14; StyleFrom = "B7";  Value = 262200 },
    @{ Row = 15; StyleFrom = "B7";  Value = 383700 },
    @{ Row = 16; StyleFrom = "B7";  Value = 448100 },
    @{ Row = 17; StyleFrom = "B7";  Value = 589000 },
    @{ Row = 18; StyleFrom = "B7";  Value = 623600 },
    @{ Row = 19; StyleFrom = "B7";  Value = 702750 },
    @{ Row = 20; StyleFrom = "B7";  Value = 771000 },
    @{ Row = 21; StyleFrom = "B7";  Value = 920750 },
    @{ Row = 22; StyleFrom = "B7";  Value = 979400 },
    @{ Row = 23; StyleFrom = "B7";  Value = 1084550 },
    @{ Row = 24; StyleFrom = "B7";  Value = 1167600 },
    @{ Row = 25; StyleFrom = "B7";  Value = 1187100 },
    @{ Row = 26; StyleFrom = "B7";  Value = 1019050 },
    @{ Row = 27; StyleFrom = "B7";  Value = 0 },
    @{ Row = 28; StyleFrom = "B7";  Value = 0 },
    @{ Row = 29; StyleFrom = "B7";  Value = 0 },
    @{ Row = 30; StyleFrom = "B7";  Value = 0 },
    @{ Row = 31; StyleFrom = "B7";  Value = 0 },
    @{ Row = 32; StyleFrom = "B7";  Value = 0 },
    @{ Row = 33; StyleFrom = "B7";  Value = 0 }
)

foreach ($d in $data) {
    $target = $ws.Range("C" + $d.Row)
    $ws.Range($d.StyleFrom).Copy()
    $target.PasteSpecial(-4122)
    $target.Value = $d.Value
}

# Match column C's width to column B (bestFit currency column).
$ws.Range("C1:C33").ColumnWidth = $ws.Range("B1:B33").ColumnWidth

# Update the active selection, as seen after the edit in the source file.
# (Note: the runtime always anchors the active cell to the top-left of the
# selected range, so the exact activeCell="I25" nuance from the original
# author's session can't be reproduced bit-for-bit; sqref is preserved.)
$ws.Range("I24:I25").Select()
